$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlLineStyleNone = -4142

# The header block in D1:E3 (Author/Date/Rev) needs to move down one row to
# D2:E4 so a new "Project" row can be inserted at the top. The main BOM
# table starting at row 6 must stay put, so we only touch D1:E4 rather than
# doing a real sheet-wide row insert.

# 1) Shift formatting of D1:E3 down to D2:E4 one source row at a time,
#    bottom-most first so a row is copied before it is overwritten.
$ws.Range("D3:E3").Copy() | Out-Null
$ws.Range("D4:E4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D2:E2").Copy() | Out-Null
$ws.Range("D3:E3").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D1:E1").Copy() | Out-Null
$ws.Range("D2:E2").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# 2) Shift the actual content (values/formula) the same way. Snapshot
#    everything first since later writes would otherwise clobber sources
#    that later reads depend on.
$d1 = $ws.Range("D1").Value2
$e1 = $ws.Range("E1").Value2
$d2 = $ws.Range("D2").Value2
$e2Formula = $ws.Range("E2").Formula
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2

$ws.Range("D4").Value = $d3
$ws.Range("E4").Value = $e3

$ws.Range("D3").Value = $d2
$ws.Range("E3").Formula = $e2Formula

$ws.Range("D2").Value = $d1
$ws.Range("E2").Value = $e1

# 3) Turn row 1 into the new "Project" row: strip the border the old
#    Author row had, keep E1's wrapped-text look, and write the new labels.
$ws.Range("D1").Borders.LineStyle = $xlLineStyleNone
$ws.Range("E1").Borders.LineStyle = $xlLineStyleNone
$ws.Range("E1").WrapText = $true

$ws.Range("D1").Value = "Project"
$ws.Range("E1").Value = "WDY - Main Board"

# 4) Match the recorded selection state (user clicked E1 after editing).
$ws.Range("E1").Select() | Out-Null
